$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'53.271.50"
$ws.Range('D2').ClearFormats()
$ws.Range('D3').Value = "'3.146.63"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +3.56%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'398.37"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.85%  '
$ws.Range('D6').Value = "'106.44"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.68%  '
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = "'0.607"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +4.42%  '
$ws.Range('D10').Value = "'38.87"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +5.73%  '
$ws.Range('E11').Value = '  +1.29%  '
$ws.Range('E12').Value = '  +1.24%  '
$ws.Range('D13').Value = "'3.645.98"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +3.65%  '
$ws.Range('D14').Value = "'18.95"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.75%  '
$ws.Range('D15').Value = "'7.96"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.03%  '
$ws.Range('E16').Value = '  +8.48%  '
$ws.Range('D17').Value = "'3.145.82"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +3.67%  '
$ws.Range('D18').Value = "'10.81"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.73%  '
$ws.Range('D19').Value = "'53.230.71"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.57%  '
$ws.Range('E20').Value = '  +4.52%  '
$ws.Range('D21').Value = "'12.89"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +4.00%  '
$ws.Range('D22').Value = "'0.0₃0976"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.36%  '
$ws.Range('D23').Value = "'71.15"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.58%  '
$ws.Range('D24').Value = "'271.17"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.60%  '
$ws.Range('D25').Value = "'3.22"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.58%  '
$ws.Range('D26').Value = "'8.04"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.60%  '
$ws.Range('D27').Value = "'27.61"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.88%  '
$ws.Range('D28').Value = "'7.50"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.63%  '
$ws.Range('E29').Value = '  +0.76%  '
$ws.Range('E31').Value = '  +2.19%  '
$ws.Range('D32').Value = "'11.02"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +7.40%  '
$ws.Range('D33').Value = "'37.24"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +6.80%  '
$ws.Range('D34').Value = "'0.0494"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +11.47%  '
$ws.Range('E35').Value = '  +0.80%  '
$ws.Range('D36').Value = "'50.45"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.51%  '
$ws.Range('D37').Value = "'3.55"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +6.73%  '
$ws.Range('D38').Value = "'1.00"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').Value = "'2.75"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +8.65%  '
$ws.Range('D40').Value = "'4.14"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +10.78%  '
$ws.Range('E41').Value = '  +1.73%  '
$ws.Range('D42').Value = "'17.37"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.40%  '
$ws.Range('E43').Value = '  +1.65%  '
$ws.Range('D44').Value = "'130.09"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +4.45%  '
$ws.Range('E45').Value = '  +1.00%  '
$ws.Range('D46').Value = "'22.31"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.49%  '
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('D49').Value = "'2.092.32"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.42%  '
$ws.Range('D50').Value = "'0.0508"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +26.94%  '
$ws.Range('D51').Value = "'0.0332"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +4.81%  '
